$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56 was a mostly-blank placeholder row (only M56 carried an empty
# quote-prefixed string). Fill it in with a new completed job-run entry,
# matching the pattern used by the surrounding rows (54/55/58...).
$ws.Range("A56").Value = "ukb51139_subset.csv"
$ws.Range("B56").Value = "28012 x 1081"
$ws.Range("C56").Value = "all"
$ws.Range("D56").Value = "no events"
$ws.Range("E56").Value = "> 160/100"
$ws.Range("F56").Value = "zscore"
$ws.Range("G56").Value = "median"
$ws.Range("H56").Value = "none"
$ws.Range("I56").Value = 25
$ws.Range("K56").Value = "N/A"
$ws.Range("L56").Value = "53.8 & 43.4"
$ws.Range("M56").Value = "-237.2 & -71.6"
$ws.Range("N56").Value = "N/A"
$ws.Range("O56").Value = "N/A"
$ws.Range("P56").Value = "25 batches"

# The blank-placeholder marker (an empty, quote-prefixed string) moves down
# to row 57's L/M cells, which were previously fully empty.
$ws.Range("L57").Value = "'"
$ws.Range("M57").Value = "'"
